# Generate Report for Handoff
#
# The "fbdf8557-9368-407b-a255-6254c559e860.md" file has finished
# translation and is now ready to be handed off. Update its status
# (and the associated handoff timestamps) across the Overview sheet
# and each per-locale detail sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is the fbdf8557... file ---
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-13-17 18:13:33"

# --- zh-cn detail sheet: row 3 is the fbdf8557... file ---
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "2016-03-17 18:13:30"

# --- de-de detail sheet: row 3 is the fbdf8557... file ---
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "2016-03-17 18:13:33"
